$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stored username value (shared string content change)
$ws.Range("A2").Value = "dilhaniwas+1@gmail.com"

# Replace the password text value with a numeric value
$ws.Range("B2").Value = 123456.0

# Shrink the custom-width column range from C:Z down to C:F,
# letting G:Z revert to the sheet's default column width
$ws.Range("G1:Z1").EntireColumn.Delete()
